$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 77388.66797673714
$ws.Range("D2").Value = 9992.97670278544
$ws.Range("E2").Value = 9770
$ws.Range("F2").Value = 19515.36047260924

$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C3").Value = 68

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 13.6
$ws.Range("H2").Value = 27.2
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 40.8
$ws.Range("K2").Value = 47.6
$ws.Range("L2").Value = 54.4
$ws.Range("M2").Value = 61.2
$ws.Range("N2").Value = 68
$ws.Range("O2").Value = 61.2
$ws.Range("P2").Value = 54.4
$ws.Range("Q2").Value = 47.6
$ws.Range("R2").Value = 34
$ws.Range("S2").Value = 20.4
$ws.Range("T2").Value = 13.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 68
$ws.Range("N3").Value = 54.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 34
$ws.Range("R3").Value = 20.4
$ws.Range("S3").Value = 13.6
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 54.4
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 10.38312417100186

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 64.3
$ws.Range("H2").Value = 53.7
$ws.Range("I2").Value = 51.6
$ws.Range("J2").Value = 1.8
$ws.Range("K2").Value = 21.6
$ws.Range("L2").Value = 33.6
$ws.Range("M2").Value = 37.8
$ws.Range("N2").Value = 42
$ws.Range("O2").Value = 30
$ws.Range("P2").Value = 25.8
$ws.Range("Q2").Value = 21.6
$ws.Range("R2").Value = 48.9
$ws.Range("S2").Value = 51.43427201306103
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.43079277624771
$ws.Range("M3").Value = 44.6
$ws.Range("N3").Value = 28.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 5.4
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 20.4
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 31
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 10.38312417100186
$ws.Range("R4").Value = 0

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("G2").Value = 183.657
$ws.Range("H2").Value = 236.82
$ws.Range("I2").Value = 287.904
$ws.Range("J2").Value = 289.686
$ws.Range("K2").Value = 311.07
$ws.Range("L2").Value = 344.334
$ws.Range("M2").Value = 381.756
$ws.Range("N2").Value = 423.336
$ws.Range("O2").Value = 453.0359999999999
$ws.Range("P2").Value = 478.578
$ws.Range("Q2").Value = 499.962
$ws.Range("R2").Value = 548.373
$ws.Range("I3").Value = 146.928
$ws.Range("J3").Value = 187.32
$ws.Range("K3").Value = 241.176
$ws.Range("L3").Value = 301.9924848484852
$ws.Range("M3").Value = 346.1464848484852
$ws.Range("N3").Value = 374.2624848484852
$ws.Range("O3").Value = 421.3864848484852
$ws.Range("P3").Value = 426.7324848484852
$ws.Range("Q3").Value = 434.6524848484852
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 146.928
$ws.Range("L4").Value = 194.052
$ws.Range("M4").Value = 224.742
$ws.Range("N4").Value = 278.598
$ws.Range("O4").Value = 325.722
$ws.Range("P4").Value = 352.65
$ws.Range("Q4").Value = 362.9292929292918

$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("P2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K4").Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 48.7
$ws.Range("J3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 4")
$ws.Range("O2").Value = 0
$ws.Range("S2").Value = 44.23427201306104
$ws.Range("T2").Value = 20.4
$ws.Range("L3").Value = 0.2307927762477035
$ws.Range("S3").Value = 9.6
$ws.Range("J4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
